# 06_finops_executivo.xlsx — apply the "docs: Atualizando README.md" edit:
#  - rename the sheet Sheet1 -> Finops_Executivo
#  - move the active selection from D20 -> B22
#  - turn A1:D16 into a real Excel Table ("Tabela1", TableStyleLight9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet rename -----------------------------------------------------
$ws.Name = "Finops_Executivo"

# --- convert the existing range into a native ListObject (Excel Table) -
# The header row (A1:D1) already carries bold/centered/bordered formatting;
# creating the table over it is what produces the table's header dxf.
$tableRange = $ws.Range("A1:D16")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Tabela1"
$lo.TableStyle = "TableStyleLight9"
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false

# --- extra dxf records for the header/table border treatment ----------
# (mirrors what Excel writes as headerRowBorderDxfId / tableBorderDxfId
# when "Format as Table" is applied over already-bordered header cells)
$headerRng = $ws.Range("A1:D1")

$fcBottom = $headerRng.FormatConditions.Add(2, "", "=TRUE")
$fcBottom.Borders.Item(1).LineStyle = 1
$fcBottom.Delete()

$fcTop = $headerRng.FormatConditions.Add(2, "", "=TRUE")
$fcTop.Borders.Item(1).LineStyle = 1
$fcTop.Delete()

# --- move the active selection -----------------------------------------
$ws.Range("B22").Select()
